$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add group2..reject headers, reusing the existing header
# style from A1 (bold, centered, bordered) via copy/paste-special of formats.
$headers = @("group1", "group2", "meandiff", "p-adj", "lower", "upper", "reject")
$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}

# Data row (row 2)
$ws.Range("A2").Value = "CSS"
$ws.Range("B2").Value = "Grassland"
$ws.Range("C2").Value = -0.5935
$ws.Range("D2").Value = 0.0232
$ws.Range("E2").Value = -1.103
$ws.Range("F2").Value = -0.08400000000000001
$ws.Range("G2").Value = $true
